$wb = $excel.ActiveWorkbook

# Sheet "展览" (展览 = sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 294
$ws1.Range("F4").Value = 1208
$ws1.Range("F5").Value = 841
$ws1.Range("F6").Value = 870
$ws1.Range("F7").Value = 1585
$ws1.Range("F11").Value = 87
$ws1.Range("F12").Value = 215
$ws1.Range("F13").Value = 66
$ws1.Range("F14").Value = 546
$ws1.Range("F15").Value = 89
$ws1.Range("F16").Value = 55
$ws1.Range("F19").Value = 305
$ws1.Range("F22").Value = 78
$ws1.Range("F23").Value = 15
$ws1.Range("F24").Value = 797

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 1052
$ws2.Range("F6").Value = 210
$ws2.Range("F8").Value = 604
$ws2.Range("F9").Value = 97

# Sheet "本地生活" (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 276

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 276
$ws4.Range("F3").Value = 294
$ws4.Range("F4").Value = 1052
$ws4.Range("F6").Value = 1208
$ws4.Range("F7").Value = 841
$ws4.Range("F8").Value = 870
$ws4.Range("F9").Value = 1585
$ws4.Range("F13").Value = 87
$ws4.Range("F14").Value = 215
$ws4.Range("F15").Value = 66
$ws4.Range("F16").Value = 546
$ws4.Range("F17").Value = 89
$ws4.Range("F18").Value = 55
$ws4.Range("F23").Value = 305
$ws4.Range("F25").Value = 210
$ws4.Range("F26").Value = 210
$ws4.Range("F29").Value = 78
$ws4.Range("F30").Value = 15
$ws4.Range("F31").Value = 797
$ws4.Range("F35").Value = 604
$ws4.Range("F36").Value = 97
$ws4.Range("F37").Value = 97

$wb.Save()
